# Replace Salary (column S) values with their natural logarithm, row by row,
# for all data rows (2 through 264) on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 264; $row++) {
    $cell = $ws.Range("S$row")
    $old = $cell.Value2
    if ($old -ne $null) {
        $cell.Value2 = $excel.Evaluate("LN($old)")
    }
}
